$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Select Sheet1 so it remains the active/tab-selected sheet
$ws.Activate()

# Add new value in cell A14 (new row at the bottom of the data)
$ws.Range("A14").Value = "changes here"

# Update selection to match the newly entered cell
$ws.Range("A14").Select()
